# Update "Chiffres COVID-19 Valais" sheet with corrected/added daily figures.
# Columns: A=Date, B=Cumul cas positifs (formula), C=Nb nouveaux cas positifs,
# D=Nb nouvelles admissions, E=Patients SI total, F=Nb intubés,
# G=Patients hospitalisés hors SI, H=Total hospitalisations (formula),
# I=Nb nouvelles sorties, J=Cumul décès (formula), K=Nb nouveaux décès (formula),
# L=Nb nouveaux décès hôpital, M=Nb nouveaux décès extra-hospitaliers.
#
# B/H/J/K are "shared" formulas that recompute automatically once the raw
# inputs (C/E/F/G/L/M) below are changed, so only the raw input cells are
# written here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small corrections to "Nb nouveaux cas positifs" (column C) on a few
#     historical rows (these ripple into the cumulative column B). ---
$ws.Cells.Item(528, 3).Value = 30
$ws.Cells.Item(533, 3).Value = 63
$ws.Cells.Item(538, 3).Value = 108
$ws.Cells.Item(545, 3).Value = 114
$ws.Cells.Item(548, 3).Value = 84
$ws.Cells.Item(575, 3).Value = 67

# --- Row 577: hospitalised-outside-ICU count correction ---
$ws.Cells.Item(577, 7).Value = 12

# --- Row 579: updated daily figures ---
$ws.Cells.Item(579, 3).Value = 11
$ws.Cells.Item(579, 5).Value = 4
$ws.Cells.Item(579, 7).Value = 13

# --- Row 580: updated daily figures ---
$ws.Cells.Item(580, 3).Value = 75
$ws.Cells.Item(580, 5).Value = 4
$ws.Cells.Item(580, 7).Value = 11

# --- Rows 581-583: newly-reported daily figures (previously blank) ---
# Columns L/M ("Nb nouveaux décès hôpital/extra-hospitaliers") carry a Text
# ("@") number format on these rows. Writing 0 straight into a Text-formatted
# cell stores it as the literal string "0"; toggle the format to General for
# the write, then restore the original Text format so the cell keeps its
# original style but holds a true numeric 0 (matches how the already-filled
# rows above were stored).
function Set-NumericOnTextCell($cell, $val) {
    $cell.NumberFormat = "General"
    $cell.Value = $val
    $cell.NumberFormat = "@"
}

$ws.Cells.Item(581, 3).Value = 34
$ws.Cells.Item(581, 5).Value = 4
$ws.Cells.Item(581, 6).Value = 2
$ws.Cells.Item(581, 7).Value = 12
Set-NumericOnTextCell $ws.Cells.Item(581, 12) 0
Set-NumericOnTextCell $ws.Cells.Item(581, 13) 0

$ws.Cells.Item(582, 3).Value = 36
$ws.Cells.Item(582, 5).Value = 4
$ws.Cells.Item(582, 6).Value = 1
$ws.Cells.Item(582, 7).Value = 13
Set-NumericOnTextCell $ws.Cells.Item(582, 12) 0
Set-NumericOnTextCell $ws.Cells.Item(582, 13) 0

$ws.Cells.Item(583, 3).Value = 2
$ws.Cells.Item(583, 5).Value = 4
$ws.Cells.Item(583, 6).Value = 1
$ws.Cells.Item(583, 7).Value = 12
Set-NumericOnTextCell $ws.Cells.Item(583, 12) 0
Set-NumericOnTextCell $ws.Cells.Item(583, 13) 0
